$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Add the newly found engineering (STC) sensor OOI barcode
$ws.Range("E34").Style = "Normal"
$ws.Range("E34").Value = "OL000374"

# Clear the bogus placeholder "RTE" calibration row (no real serial/barcode info)
$ws.Range("A35:D35").ClearContents()
$ws.Range("F35").ClearContents()
$ws.Range("I35").ClearContents()

# Stray formatting left over from editing near column L (left align an
# unused range, matching the workbook's on-disk diff)
$ws.Range("L3:L10").HorizontalAlignment = -4131
